# Contacts Final - 14 Oct 2024
# Insert a new "Country" column between the existing "MailingState" (K) and
# "Status" (old L, now M) columns on the Contact sheet, shifting the
# remaining columns (Status..LineOfBusiness) one position to the right,
# and populate it with "United States" for both data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")
$ws.Activate()

# Insert a whole new column at L; this shifts L:P -> M:Q and keeps
# header styling/number formats intact automatically.
$ws.Columns("L").Insert()

# New header + values for the inserted "Country" column.
$ws.Range("L1").Value = "Country"
$ws.Range("L2").Value = "United States"
$ws.Range("L3").Value = "United States"

# Match the column width used for the new column in the final workbook.
$ws.Columns("L").ColumnWidth = 13.166666666666666

# Reproduce the updated view state (top-left cell / selection) as closely
# as possible.
$ws.Range("L2:L3").Select()
